{"js": "// Replace the 100 arithmetic-equation answers in the table, in document\n// order, per the commit's diff. The document body is: 1 heading paragraph\n// (\"2024-08-12 Monday\") followed by 100 table-cell paragraphs (one run\n// each) holding the equations.\nconst replacements = [\n  [\"90-72=18\", \"32+9=41\"],\n  [\"50-31=19\", \"26+27=53\"],\n  [\"74+8=82\", \"24+49=73\"],\n  [\"17+64=81\", \"47-39=8\"],\n  [\"39+47=86\", \"8+9=17\"],\n  [\"95-17=78\", \"39+38=77\"],\n  [\"96-49=47\", \"85-49=36\"],\n  [\"25+29=54\", \"42-9=33\"],\n  [\"43-17=26\", \"63-17=46\"],\n  [\"67+27=94\", \"90-38=52\"],\n  [\"9+46=55\", \"88-79=9\"],\n  [\"38+44=82\", \"25+6=31\"],\n  [\"98-79=19\", \"37+45=82\"],\n  [\"65-38=27\", \"63-57=6\"],\n  [\"54-46=8\", \"21-4=17\"],\n  [\"13-9=4\", \"95-48=47\"],\n  [\"6+48=54\", \"80-35=45\"],\n  [\"92-67=25\", \"81-39=42\"],\n  [\"82+9=91\", \"9+13=22\"],\n  [\"87+5=92\", \"72-26=46\"],\n  [\"83-66=17\", \"53-39=14\"],\n  [\"30-28=2\", \"16+57=73\"],\n  [\"88+9=97\", \"78+15=93\"],\n  [\"81-52=29\", \"6+17=23\"],\n  [\"93-4=89\", \"92-27=65\"],\n  [\"61-32=29\", \"69+25=94\"],\n  [\"92-57=35\", \"78+16=94\"],\n  [\"83-24=59\", \"20-12=8\"],\n  [\"91-75=16\", \"62-4=58\"],\n  [\"45+8=53\", \"41-14=27\"],\n  [\"82-49=33\", \"82-35=47\"],\n  [\"5+47=52\", \"8+88=96\"],\n  [\"56+19=75\", \"61-54=7\"],\n  [\"17+25=42\", \"72-44=28\"],\n  [\"8+86=94\", \"73-17=56\"],\n  [\"7+15=22\", \"90-44=46\"],\n  [\"41-4=37\", \"55-26=29\"],\n  [\"30-7=23\", \"85-46=39\"],\n  [\"94-38=56\", \"39+54=93\"],\n  [\"50-38=12\", \"38+3=41\"],\n  [\"61-12=49\", \"59+14=73\"],\n  [\"58+23=81\", \"41-18=23\"],\n  [\"62-3=59\", \"33+9=42\"],\n  [\"9+84=93\", \"48+6=54\"],\n  [\"36-29=7\", \"62-49=13\"],\n  [\"75-67=8\", \"76+17=93\"],\n  [\"36+48=84\", \"28+7=35\"],\n  [\"77-39=38\", \"26+66=92\"],\n  [\"5+7=12\", \"60-31=29\"],\n  [\"90-12=78\", \"47+9=56\"],\n  [\"66+19=85\", \"18+66=84\"],\n  [\"9+62=71\", \"9+77=86\"],\n  [\"41-36=5\", \"47+24=71\"],\n  [\"67+14=81\", \"44+47=91\"],\n  [\"19+3=22\", \"63-59=4\"],\n  [\"46+15=61\", \"27+45=72\"],\n  [\"47+36=83\", \"59+36=95\"],\n  [\"43-37=6\", \"29+59=88\"],\n  [\"15+8=23\", \"29+57=86\"],\n  [\"51-8=43\", \"63-7=56\"],\n  [\"65-59=6\", \"94-17=77\"],\n  [\"93-89=4\", \"51-17=34\"],\n  [\"72-55=17\", \"62-25=37\"],\n  [\"20-4=16\", \"22+29=51\"],\n  [\"4+9=13\", \"28+18=46\"],\n  [\"85-9=76\", \"37+17=54\"],\n  [\"95-57=38\", \"39+28=67\"],\n  [\"59+39=98\", \"20-2=18\"],\n  [\"96-8=88\", \"90-68=22\"],\n  [\"46+39=85\", \"42+49=91\"],\n  [\"13-5=8\", \"60-36=24\"],\n  [\"28+56=84\", \"54-49=5\"],\n  [\"91-19=72\", \"65+28=93\"],\n  [\"8+36=44\", \"28+25=53\"],\n  [\"2+89=91\", \"49+8=57\"],\n  [\"28+6=34\", \"2+29=31\"],\n  [\"69+24=93\", \"8+34=42\"],\n  [\"21-16=5\", \"53-24=29\"],\n  [\"52+9=61\", \"18+18=36\"],\n  [\"25+48=73\", \"27+4=31\"],\n  [\"70-21=49\", \"52-8=44\"],\n  [\"83-18=65\", \"28+46=74\"],\n  [\"95-28=67\", \"41-15=26\"],\n  [\"49+35=84\", \"74-39=35\"],\n  [\"45-36=9\", \"92-63=29\"],\n  [\"37+59=96\", \"51-17=34\"],\n  [\"37+27=64\", \"83-28=55\"],\n  [\"63-37=26\", \"7+45=52\"],\n  [\"59+13=72\", \"29+35=64\"],\n  [\"84-36=48\", \"57+19=76\"],\n  [\"96-48=48\", \"18+9=27\"],\n  [\"84-68=16\", \"60-42=18\"],\n  [\"84-78=6\", \"19+16=35\"],\n  [\"23+59=82\", \"14+68=82\"],\n  [\"24+8=32\", \"60-21=39\"],\n  [\"69+23=92\", \"69+9=78\"],\n  [\"72-45=27\", \"31-8=23\"],\n  [\"9+25=34\", \"50-34=16\"],\n  [\"74-45=29\", \"72-57=15\"],\n  [\"73-17=56\", \"37+14=51\"]\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Skip the heading paragraph; the remaining paragraphs are the table cells\n// in row-major, left-to-right order - the same order the diff lists them.\nconst cellParagraphs = paragraphs.items.slice(1);\n\nif (cellParagraphs.length !== replacements.length) {\n  throw new Error(\n    \"Expected \" + replacements.length + \" table-cell paragraphs, found \" +\n    cellParagraphs.length\n  );\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = cellParagraphs[i];\n  const current = para.text.trim();\n  if (current !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected \\\"\" + oldText +\n      \"\\\" but found \\\"\" + current + \"\\\"\"\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-equation answers in the table, in document\n# order (row-major, left-to-right), per the commit's diff. The table is the\n# sole table in the document: 20 rows x 5 columns, one run per cell.\n$replacements = @(\n    @(\"90-72=18\", \"32+9=41\"),\n    @(\"50-31=19\", \"26+27=53\"),\n    @(\"74+8=82\", \"24+49=73\"),\n    @(\"17+64=81\", \"47-39=8\"),\n    @(\"39+47=86\", \"8+9=17\"),\n    @(\"95-17=78\", \"39+38=77\"),\n    @(\"96-49=47\", \"85-49=36\"),\n    @(\"25+29=54\", \"42-9=33\"),\n    @(\"43-17=26\", \"63-17=46\"),\n    @(\"67+27=94\", \"90-38=52\"),\n    @(\"9+46=55\", \"88-79=9\"),\n    @(\"38+44=82\", \"25+6=31\"),\n    @(\"98-79=19\", \"37+45=82\"),\n    @(\"65-38=27\", \"63-57=6\"),\n    @(\"54-46=8\", \"21-4=17\"),\n    @(\"13-9=4\", \"95-48=47\"),\n    @(\"6+48=54\", \"80-35=45\"),\n    @(\"92-67=25\", \"81-39=42\"),\n    @(\"82+9=91\", \"9+13=22\"),\n    @(\"87+5=92\", \"72-26=46\"),\n    @(\"83-66=17\", \"53-39=14\"),\n    @(\"30-28=2\", \"16+57=73\"),\n    @(\"88+9=97\", \"78+15=93\"),\n    @(\"81-52=29\", \"6+17=23\"),\n    @(\"93-4=89\", \"92-27=65\"),\n    @(\"61-32=29\", \"69+25=94\"),\n    @(\"92-57=35\", \"78+16=94\"),\n    @(\"83-24=59\", \"20-12=8\"),\n    @(\"91-75=16\", \"62-4=58\"),\n    @(\"45+8=53\", \"41-14=27\"),\n    @(\"82-49=33\", \"82-35=47\"),\n    @(\"5+47=52\", \"8+88=96\"),\n    @(\"56+19=75\", \"61-54=7\"),\n    @(\"17+25=42\", \"72-44=28\"),\n    @(\"8+86=94\", \"73-17=56\"),\n    @(\"7+15=22\", \"90-44=46\"),\n    @(\"41-4=37\", \"55-26=29\"),\n    @(\"30-7=23\", \"85-46=39\"),\n    @(\"94-38=56\", \"39+54=93\"),\n    @(\"50-38=12\", \"38+3=41\"),\n    @(\"61-12=49\", \"59+14=73\"),\n    @(\"58+23=81\", \"41-18=23\"),\n    @(\"62-3=59\", \"33+9=42\"),\n    @(\"9+84=93\", \"48+6=54\"),\n    @(\"36-29=7\", \"62-49=13\"),\n    @(\"75-67=8\", \"76+17=93\"),\n    @(\"36+48=84\", \"28+7=35\"),\n    @(\"77-39=38\", \"26+66=92\"),\n    @(\"5+7=12\", \"60-31=29\"),\n    @(\"90-12=78\", \"47+9=56\"),\n    @(\"66+19=85\", \"18+66=84\"),\n    @(\"9+62=71\", \"9+77=86\"),\n    @(\"41-36=5\", \"47+24=71\"),\n    @(\"67+14=81\", \"44+47=91\"),\n    @(\"19+3=22\", \"63-59=4\"),\n    @(\"46+15=61\", \"27+45=72\"),\n    @(\"47+36=83\", \"59+36=95\"),\n    @(\"43-37=6\", \"29+59=88\"),\n    @(\"15+8=23\", \"29+57=86\"),\n    @(\"51-8=43\", \"63-7=56\"),\n    @(\"65-59=6\", \"94-17=77\"),\n    @(\"93-89=4\", \"51-17=34\"),\n    @(\"72-55=17\", \"62-25=37\"),\n    @(\"20-4=16\", \"22+29=51\"),\n    @(\"4+9=13\", \"28+18=46\"),\n    @(\"85-9=76\", \"37+17=54\"),\n    @(\"95-57=38\", \"39+28=67\"),\n    @(\"59+39=98\", \"20-2=18\"),\n    @(\"96-8=88\", \"90-68=22\"),\n    @(\"46+39=85\", \"42+49=91\"),\n    @(\"13-5=8\", \"60-36=24\"),\n    @(\"28+56=84\", \"54-49=5\"),\n    @(\"91-19=72\", \"65+28=93\"),\n    @(\"8+36=44\", \"28+25=53\"),\n    @(\"2+89=91\", \"49+8=57\"),\n    @(\"28+6=34\", \"2+29=31\"),\n    @(\"69+24=93\", \"8+34=42\"),\n    @(\"21-16=5\", \"53-24=29\"),\n    @(\"52+9=61\", \"18+18=36\"),\n    @(\"25+48=73\", \"27+4=31\"),\n    @(\"70-21=49\", \"52-8=44\"),\n    @(\"83-18=65\", \"28+46=74\"),\n    @(\"95-28=67\", \"41-15=26\"),\n    @(\"49+35=84\", \"74-39=35\"),\n    @(\"45-36=9\", \"92-63=29\"),\n    @(\"37+59=96\", \"51-17=34\"),\n    @(\"37+27=64\", \"83-28=55\"),\n    @(\"63-37=26\", \"7+45=52\"),\n    @(\"59+13=72\", \"29+35=64\"),\n    @(\"84-36=48\", \"57+19=76\"),\n    @(\"96-48=48\", \"18+9=27\"),\n    @(\"84-68=16\", \"60-42=18\"),\n    @(\"84-78=6\", \"19+16=35\"),\n    @(\"23+59=82\", \"14+68=82\"),\n    @(\"24+8=32\", \"60-21=39\"),\n    @(\"69+23=92\", \"69+9=78\"),\n    @(\"72-45=27\", \"31-8=23\"),\n    @(\"9+25=34\", \"50-34=16\"),\n    @(\"74-45=29\", \"72-57=15\"),\n    @(\"73-17=56\", \"37+14=51\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nif (($t.Rows.Count * $t.Columns.Count) -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) cells, found $($t.Rows.Count * $t.Columns.Count)\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($cI = 1; $cI -le $t.Columns.Count; $cI++) {\n        $pair = $replacements[$idx]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n\n        $cell = $t.Cell($r, $cI)\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $oldText) {\n            throw \"Cell ($r,$cI) text mismatch: expected '$oldText' but found '$current'\"\n        }\n        $cell.Range.Text = $newText\n\n        $idx++\n    }\n}\n\nWrite-Output \"Replaced $idx cells\"\n"}
